# Apply the Team_PER_2004 update:
#  - rows are re-sorted into a new team order (the underlying Python/pandas
#    rebuild changed row order, which is why the shared-string table order
#    changes too -- Excel just records strings in the order it first meets
#    them while writing rows top-to-bottom)
#  - the PER column (C) is fixed from a season-total figure to a per-game
#    average (old total / games played)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order (row number => Team code), matches the reordered shared
# strings table in the target workbook.
$teams = @{
    2  = "POR"
    3  = "NJN"
    4  = "CLE"
    5  = "DAL"
    6  = "MIA"
    7  = "SEA"
    8  = "ATL"
    9  = "WAS"
    10 = "MIL"
    11 = "LAC"
    12 = "SAS"
    13 = "DET"
    14 = "ORL"
    15 = "UTA"
    16 = "MEM"
    17 = "HOU"
    18 = "DEN"
    19 = "LAL"
    20 = "GSW"
    21 = "IND"
    22 = "CHI"
    23 = "PHI"
    24 = "BOS"
    25 = "TOR"
    26 = "SAC"
    27 = "PHO"
    28 = "NOH"
    29 = "NYK"
    30 = "MIN"
}

# Corrected, per-game PER values for the same rows.
$per = @{
    2  = 10.18
    3  = 11.3
    4  = 11.04
    5  = 16.65833333333333
    6  = 10.65625
    7  = 15.51428571428571
    8  = 11.06
    9  = 13.26923076923077
    10 = 15.07272727272727
    11 = 10.97333333333333
    12 = 11.92666666666667
    13 = 12.94545454545454
    14 = 11.26153846153846
    15 = 12.11428571428571
    16 = 11.75714285714286
    17 = 12.26428571428571
    18 = 11.82857142857143
    19 = 13.92
    20 = 12.59375
    21 = 13.13333333333333
    22 = 11.8
    23 = 12.9
    24 = 14.75
    25 = 10.88888888888889
    26 = 14.56153846153846
    27 = 13.4625
    28 = 11.9
    29 = 12.7125
    30 = 12.4375
}

for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 2).Value = $teams[$r]
    $ws.Cells.Item($r, 3).Value = $per[$r]
}
